# Update odds on existing rows 3-6 (FlashScore weekly refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (River Plate vs Miramar) ---
$ws.Range("G3").Value = 2.45
$ws.Range("I3").Value = 3.1
$ws.Range("J3").Value = 3.1
$ws.Range("L3").Value = 3.6
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 1.7
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.83
$ws.Range("W3").Value = 7.5
$ws.Range("X3").Value = 11
$ws.Range("AI3").Value = 15
$ws.Range("AJ3").Value = 12
$ws.Range("AL3").Value = 26
$ws.Range("AN3").Value = 4.33
$ws.Range("AO3").Value = 13
$ws.Range("AQ3").Value = 41

# --- Row 4 (Fenix vs Nacional) ---
$ws.Range("G4").Value = 7.5
$ws.Range("H4").Value = 4.75
$ws.Range("I4").Value = 1.42
$ws.Range("J4").Value = 7.5
$ws.Range("K4").Value = 2.3
$ws.Range("L4").Value = 1.95
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 3.4
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.88
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("AD4").Value = 9
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 81
$ws.Range("AI4").Value = 6
$ws.Range("AK4").Value = 9
$ws.Range("AM4").Value = 34
$ws.Range("AN4").Value = 8.5
$ws.Range("AR4").Value = 201
$ws.Range("AU4").Value = 10

# --- Row 5 (Cerro Largo vs Liverpool M.) ---
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75

# --- Row 6 (Defensor Sp. vs Progreso) ---
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.85

# --- Insert a new match (Carabobo vs Monagas) before the existing row 7,
#     pushing the current row 7 (La Guaira vs Estudiantes Merida) down to row 8 ---
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "Qa7iAtsI"
$ws.Range("B7").Value = "13/11/2024"
$ws.Range("C7").Value = "19:30"
$ws.Range("D7").Value = "VENEZUELA - LIGA FUTVE"
$ws.Range("E7").Value = "Carabobo"
$ws.Range("F7").Value = "Monagas"
$ws.Range("G7").Value = 1.98
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 3.55
$ws.Range("J7").Value = 2.62
$ws.Range("K7").Value = 2.02
$ws.Range("L7").Value = 4.1
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 6.3
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 2.65
$ws.Range("Q7").Value = 2.05
$ws.Range("R7").Value = 1.6
$ws.Range("S7").Value = 1.45
$ws.Range("T7").Value = 2.37
$ws.Range("U7").Value = 1.9
$ws.Range("V7").Value = 1.72
$ws.Range("W7").Value = 6.3
$ws.Range("X7").Value = 8.5
$ws.Range("Y7").Value = 8.75
$ws.Range("Z7").Value = 16.5
$ws.Range("AA7").Value = 17.5
$ws.Range("AB7").Value = 35
$ws.Range("AC7").Value = 8.5
$ws.Range("AD7").Value = 6.5
$ws.Range("AE7").Value = 17
$ws.Range("AF7").Value = 100
$ws.Range("AG7").Value = 900
$ws.Range("AH7").Value = 9
$ws.Range("AI7").Value = 17.5
$ws.Range("AJ7").Value = 12.5
$ws.Range("AK7").Value = 50
$ws.Range("AL7").Value = 35
$ws.Range("AM7").Value = 50
$ws.Range("AN7").Value = 3.7
$ws.Range("AO7").Value = 10.25
$ws.Range("AP7").Value = 22
$ws.Range("AQ7").Value = 40
$ws.Range("AR7").Value = 90
$ws.Range("AS7").Value = 350
$ws.Range("AT7").Value = 2.35
$ws.Range("AU7").Value = 7.9
$ws.Range("AV7").Value = 90
$ws.Range("AW7").Value = 5.2
$ws.Range("AX7").Value = 20
$ws.Range("AY7").Value = 32
$ws.Range("AZ7").Value = 110
$ws.Range("BA7").Value = 175
$ws.Range("BB7").Value = 500
$ws.Range("BC7").Value = 51
$ws.Range("BD7").Value = 51

# --- Row 8 (now La Guaira vs Estudiantes Merida, shifted down) odds refresh ---
$ws.Range("J8").Value = 2.62
$ws.Range("K8").Value = 2.07
$ws.Range("L8").Value = 3.85
$ws.Range("N8").Value = 7.6
$ws.Range("P8").Value = 3.2
$ws.Range("S8").Value = 1.37
$ws.Range("T8").Value = 2.85
$ws.Range("W8").Value = 8.5
$ws.Range("X8").Value = 11
$ws.Range("AA8").Value = 15.5
$ws.Range("AB8").Value = 23
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 6.4
$ws.Range("AH8").Value = 10.5
$ws.Range("AI8").Value = 18.5
$ws.Range("AJ8").Value = 11.5
$ws.Range("AL8").Value = 29
$ws.Range("AO8").Value = 10.5
$ws.Range("AP8").Value = 17.5
$ws.Range("AR8").Value = 65
$ws.Range("AX8").Value = 18.5
$ws.Range("AY8").Value = 25
$ws.Range("AZ8").Value = 100
$ws.Range("BA8").Value = 120
$ws.Range("BB8").Value = 300
